# Add columns I (I0) and J (IF) to the sheet, matching the style/format
# already used for the other header/data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1), copy formatting from the existing header cell H1 so the
# new header cells end up bold/centered/bordered like the rest of the header
# row (same cell style the other headers already use).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-38.
$values = @{
    2  = @(8, 8)
    3  = @(5, 7)
    4  = @(9, 9)
    5  = @(7, 7)
    6  = @(9, 9)
    7  = @(9, 9)
    8  = @(7, 7)
    9  = @(6, 7)
    10 = @(6, 6)
    11 = @(8, 8)
    12 = @(8, 9)
    13 = @(8, 8)
    14 = @(6, 7)
    15 = @(8, 9)
    16 = @(7, 7)
    17 = @(7, 8)
    18 = @(4, 6)
    19 = @(6, 7)
    20 = @(8, 8)
    21 = @(5, 6)
    22 = @(1, 4)
    23 = @(1, 2)
    24 = @(1, 5)
    25 = @(1, 4)
    26 = @(1, 4)
    27 = @(1, 3)
    28 = @(1, 4)
    29 = @(1, 5)
    30 = @(1, 6)
    31 = @(1, 5)
    32 = @(1, 5)
    33 = @(1, 1)
    34 = @(1, 5)
    35 = @(1, 5)
    36 = @(1, 3)
    37 = @(4, 6)
    38 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
